$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Copy the header-row style from L4 onto the two new header cells
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Match the author's final selection: M4, anchored over M4:N4
$ws.Range("M4:N4").Select() | Out-Null
